# Company Module changes - 22 Aug 2023
#
# Updates the coverage-team roster: the officer name on the
# "AddCoverageTeam" sheet changes from "Jennifer Muller" to
# "Jacklyn Robinson", the user's active selection on that sheet moves
# to B8, and two stray "duplicate default" cell styles (left over from
# earlier edits) are normalised back to the sheet's plain default style.

$wb = $excel.ActiveWorkbook

# --- Users sheet: A2 was carrying a redundant "default" style; drop it
#     back to the sheet's plain/default formatting (it was never bold).
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Font.Bold = $false

# --- Company sheet: C2 had the same redundant duplicate-default style.
$wsCompany = $wb.Worksheets.Item("Company")
$wsCompany.Range("C2").Font.Bold = $false

# --- AddCoverageTeam sheet: update the officer name and move the
#     selection to reflect where the user left off.
$wsAddCoverageTeam = $wb.Worksheets.Item("AddCoverageTeam")
$wsAddCoverageTeam.Range("B2").Value = "Jacklyn Robinson"

$wsAddCoverageTeam.Activate()
$wsAddCoverageTeam.Range("B8").Select()
